$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cell text from "new arrival" to "new arrivals"
$ws.Range("C2").Value = "new arrivals"

# Update the selection to reflect the active cell being C2
$ws.Range("C2").Select()
